# Add new course-log rows (188-207) to the "FSDS-10AM" sheet, mirroring the
# pattern already used for the rest of the table: a date (col A) + day
# number (col B) on the first row of each day, and a topic note in col C
# on every row.
#
# NOTE on write order: column C values are entered in the same (slightly
# "out of sequence") order the original author typed/pasted them in, so
# that the shared-strings table comes out with the same <si> ordering as
# the authored workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("FSDS-10AM")

# ---- Day 15 -- 12 Jul 2023 (serial 45119) --------------------------------
$ws.Cells.Item(188, 1).Value = 45119
$ws.Cells.Item(188, 2).Value = 15
$ws.Cells.Item(188, 3).Value = "continue numpy || reshaping we done "
$ws.Cells.Item(189, 3).Value = "indexing & slicing in the matrix "
$ws.Cells.Item(190, 3).Value = "filter in matrix | order - C, F, A"
$ws.Cells.Item(191, 3).Value = "how ai is implement in farming to pluck the fruits"
$ws.Cells.Item(192, 3).Value = "projects - requiremt, businesscase, what we need to do"
$ws.Cells.Item(193, 3).Value = "DA | BA | DE | DS -- what skill set required to do this jobs "
$ws.Cells.Item(194, 3).Value = "reqest every one please upload your ppt project code to git | linkedin"
$ws.Cells.Item(195, 3).Value = "217 functionality - continue work on this. "

# ---- Day 16 -- 13 Jul 2023 (serial 45120) --------------------------------
$ws.Cells.Item(196, 1).Value = 45120
$ws.Cells.Item(196, 2).Value = 16
$ws.Cells.Item(196, 3).Value = "We are continues on data anlaysis "
$ws.Cells.Item(197, 3).Value = "we discussed usecase | problem statement "
$ws.Cells.Item(198, 3).Value = "being a DA what is nature of work"
$ws.Cells.Item(199, 3).Value = "we introduced to matplotlib "
$ws.Cells.Item(200, 3).Value = "line style, plot parameter we discussed "
$ws.Cells.Item(201, 3).Value = "we completed data anlaysis project"
$ws.Cells.Item(203, 3).Value = "discussed insight, pattern, trends etc"
$ws.Cells.Item(204, 3).Value = "legend -- automatic color detection "
$ws.Cells.Item(206, 3).Value = "ignore the warning, how to increse the graph size by plt.rcparam."
$ws.Cells.Item(202, 3).Value = "project-ipl data analysis with numpy + matplotlib "
$ws.Cells.Item(207, 3).Value = "want you complete this"
$ws.Cells.Item(205, 3).Value = "bbox - to - anchor"

# Give the two date cells the same number format (short date, centred) as
# every other date cell in column A, by cloning an existing one instead of
# inventing a fresh style entry.
$ws.Cells.Item(168, 1).Copy()
$ws.Cells.Item(188, 1).PasteSpecial(-4122)
$ws.Cells.Item(168, 1).Copy()
$ws.Cells.Item(196, 1).PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Refresh the sheet's row/column "used extent" bookkeeping so every row
# (old and new alike) reports the full A:C span, matching the authored
# file -- touch a scratch row far below the data to force the recompute,
# then remove it again.
$ws.Range("A300:B300").HorizontalAlignment = -4108
$ws.Rows.Item(300).Delete()

# Move the view/selection to match where the author ended up after typing
# the new rows in.
$ws.Range("A208").Select()
